$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.196.95'
$ws.Range("E2").Value = '  +0.22%  '
$ws.Range("D3").Value = '2.366.93'
$ws.Range("E3").Value = '  +0.70%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '548.35'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +0.64%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.46'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  -1.05%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +5.27%  '
$ws.Range("E9").Value = '  +3.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.57'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  +2.74%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.354'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  -1.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.14'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = '  +2.39%  '
$ws.Range("D14").Value = '2.789.46'
$ws.Range("E14").Value = '  +0.77%  '
$ws.Range("D15").Value = '58.111.21'
$ws.Range("E15").Value = '  +0.12%  '
$ws.Range("E16").Value = '  +1.90%  '
$ws.Range("D17").Value = '2.366.25'
$ws.Range("E17").Value = '  +0.78%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.99'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  +3.44%  '
$ws.Range("E19").Value = '  +2.68%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '330.85'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  -1.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.89'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +2.90%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '63.49'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  +2.84%  '
$ws.Range("E24").Value = '  -1.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  +0.19%  '
$ws.Range("E26").Value = '  -2.36%  '
$ws.Range("E27").Value = '  -6.13%  '
$ws.Range("E28").Value = '  -0.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '170.25'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  -0.12%  '
$ws.Range("D30").Value = '0.0₃0742'
$ws.Range("E30").Value = '  +1.54%  '
$ws.Range("E31").Value = '  +0.38%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.43'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  -0.23%  '
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.996'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  -4.15%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.18'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  -0.02%  '
$ws.Range("E37").Value = '  -1.47%  '
$ws.Range("E38").Value = '  -2.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.412'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  +8.82%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '143.12'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  -3.60%  '
$ws.Range("E41").Value = '  +2.32%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '288.34'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  +0.52%  '
$ws.Range("E43").Value = '  +2.55%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0517'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  +2.44%  '
$ws.Range("E45").Value = '  +0.54%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '18.89'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  -2.13%  '
$ws.Range("E47").Value = '  +2.53%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.390'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  +2.26%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '11.08'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  +0.21%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.72'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  +0.78%  '
$ws.Range("E51").Value = '  +0.12%  '
